$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to stay as text so Excel does not auto-convert
# numeric-looking strings (e.g. "1.003") into floating point numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.126.19"
$ws.Range("E2").Value = "  -0.23%  "
$ws.Range("D3").Value = "1.669.42"
$ws.Range("E3").Value = "  -0.62%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("D5").Value = "210.47"
$ws.Range("E5").Value = "  -2.49%  "
$ws.Range("D6").Value = "0.5223"
$ws.Range("E6").Value = "  -0.79%  "
$ws.Range("E7").Value = "  -0.35%  "
$ws.Range("D8").Value = "0.2622"
$ws.Range("E8").Value = "  -2.31%  "
$ws.Range("D9").Value = "0.06322"
$ws.Range("E9").Value = "  -0.56%  "
$ws.Range("E10").Value = "  -0.89%  "
$ws.Range("D11").Value = "0.07555"
$ws.Range("E11").Value = "  -0.87%  "
$ws.Range("D12").Value = "1.680.23"
$ws.Range("E12").Value = "  -0.10%  "
$ws.Range("D13").Value = "4.428"
$ws.Range("E13").Value = "  -2.09%  "
$ws.Range("D14").Value = "0.5477"
$ws.Range("E14").Value = "  -4.34%  "
$ws.Range("D15").Value = "0.000008023"
$ws.Range("E15").Value = "  -2.25%  "
$ws.Range("D16").Value = "66.45"
$ws.Range("D17").Value = "26.165.87"
$ws.Range("E17").Value = "  -0.23%  "
$ws.Range("E18").Value = "  -0.34%  "
$ws.Range("D19").Value = "4.750"
$ws.Range("E19").Value = "  -2.35%  "
$ws.Range("D20").Value = "187.53"
$ws.Range("E20").Value = "  -1.19%  "
$ws.Range("E21").Value = "  -3.88%  "
$ws.Range("D22").Value = "6.239"
$ws.Range("E22").Value = "  +0.28%  "
$ws.Range("D23").Value = "1.004"
$ws.Range("E23").Value = "  -0.30%  "
$ws.Range("D24").Value = "149.55"
$ws.Range("E24").Value = "  +0.37%  "
$ws.Range("E25").Value = "  -1.11%  "
$ws.Range("D26").Value = "7.480"
$ws.Range("E26").Value = "  -2.94%  "
$ws.Range("D27").Value = "15.79"
$ws.Range("E27").Value = "  -0.80%  "
$ws.Range("D28").Value = "0.06301"
$ws.Range("E28").Value = "  -1.82%  "
$ws.Range("D29").Value = "1.358"
$ws.Range("E29").Value = "  -0.49%  "
$ws.Range("D30").Value = "1.282"
$ws.Range("E30").Value = "  -2.32%  "
$ws.Range("D31").Value = "3.526"
$ws.Range("E31").Value = "  -0.79%  "
$ws.Range("D32").Value = "3.420"
$ws.Range("E32").Value = "  -4.01%  "
$ws.Range("E33").Value = "  -1.45%  "
$ws.Range("D34").Value = "1.005"
$ws.Range("E34").Value = "  -1.02%  "
$ws.Range("D35").Value = "0.6017"
$ws.Range("E35").Value = "  -1.05%  "
$ws.Range("D36").Value = "2.769"
$ws.Range("E36").Value = "  +0.91%  "
$ws.Range("E37").Value = "  -1.08%  "
$ws.Range("D38").Value = "1.116.23"
$ws.Range("E38").Value = "  +2.18%  "
$ws.Range("D39").Value = "0.01613"
$ws.Range("E39").Value = "  -1.52%  "
$ws.Range("D40").Value = "6.067"
$ws.Range("E40").Value = "  -1.52%  "
$ws.Range("D41").Value = "0.8637"
$ws.Range("E41").Value = "  -1.66%  "
$ws.Range("E42").Value = "  -0.64%  "
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("D44").Value = "1.824.13"
$ws.Range("E44").Value = "  -0.43%  "
$ws.Range("D45").Value = "0.00000000111"
$ws.Range("E45").Value = "  +0.79%  "
$ws.Range("D46").Value = "55.51"
$ws.Range("E46").Value = "  -3.34%  "
$ws.Range("D47").Value = "1.001"
$ws.Range("E47").Value = "  +0.20%  "
$ws.Range("D48").Value = "8.081"
$ws.Range("E48").Value = "  -0.03%  "
$ws.Range("D49").Value = "0.05237"
$ws.Range("E49").Value = "  -0.42%  "
$ws.Range("D50").Value = "0.4238"
$ws.Range("E50").Value = "  -0.95%  "
$ws.Range("D51").Value = "5.920"
$ws.Range("E51").Value = "  -1.40%  "

# Remove the temporary text number-format again so the cell style
# matches the original (no explicit style index on data cells).
$ws.Range("D2:D51").ClearFormats()

